$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.38"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.88"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.386"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06008"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.389"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8171"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9523"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1433"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07434"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03290"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03047"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09411"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.005"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001592"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04802"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005902"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005472"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.004155"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009871"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.673"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.421"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.00007003"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03995"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006473"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1073"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002901"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005828"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005129"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.004164"
